$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.720.73"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").Value = "1.901.61"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'311.68"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").Value = "'0.5222"
$ws.Range("E7").Value = "  +5.72%  "
$ws.Range("D8").Value = "'0.3784"
$ws.Range("E8").Value = "  -0.25%  "
$ws.Range("D9").Value = "'0.07238"
$ws.Range("E9").Value = "  -1.22%  "
$ws.Range("D10").Value = "'21.26"
$ws.Range("E10").Value = "  +3.26%  "
$ws.Range("E11").Value = "  -0.89%  "
$ws.Range("E12").Value = "  +0.33%  "
$ws.Range("D13").Value = "1.922.51"
$ws.Range("E13").Value = "  +3.86%  "
$ws.Range("D14").Value = "'5.446"
$ws.Range("E14").Value = "  -0.24%  "
$ws.Range("D15").Value = "'92.36"
$ws.Range("E15").Value = "  +1.49%  "
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("D17").Value = "'0.000008682"
$ws.Range("E17").Value = "  -0.55%  "
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("D19").Value = "27.754.90"
$ws.Range("E19").Value = "  +0.38%  "
$ws.Range("E20").Value = "  -0.29%  "
$ws.Range("E21").Value = "  +0.46%  "
$ws.Range("D22").Value = "2.160.96"
$ws.Range("E22").Value = "  +0.43%  "
$ws.Range("E23").Value = "  +0.70%  "
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("D25").Value = "'152.97"
$ws.Range("E25").Value = "  -0.69%  "
$ws.Range("D26").Value = "'1.864"
$ws.Range("E26").Value = "  +0.84%  "
$ws.Range("E27").Value = "  -0.50%  "
$ws.Range("D28").Value = "'2.160"
$ws.Range("E28").Value = "  -1.03%  "
$ws.Range("D29").Value = "'114.43"
$ws.Range("E29").Value = "  -0.56%  "
$ws.Range("D30").Value = "'4.839"
$ws.Range("E30").Value = "  -0.76%  "
$ws.Range("D31").Value = "'0.09082"
$ws.Range("E31").Value = "  +1.67%  "
$ws.Range("D32").Value = "'3.188"
$ws.Range("E32").Value = "  -2.09%  "
$ws.Range("D33").Value = "'4.840"
$ws.Range("E33").Value = "  +4.31%  "
$ws.Range("D34").Value = "'0.7776"
$ws.Range("E34").Value = "  +1.50%  "
$ws.Range("D35").Value = "'1.219"
$ws.Range("E35").Value = "  -0.72%  "
$ws.Range("D36").Value = "'0.02092"
$ws.Range("E36").Value = "  +2.11%  "
$ws.Range("D37").Value = "'2.574"
$ws.Range("E37").Value = "  +0.77%  "
$ws.Range("D38").Value = "'3.073"
$ws.Range("E38").Value = "  +2.80%  "
$ws.Range("E39").Value = "  -0.53%  "
$ws.Range("D40").Value = "'0.5539"
$ws.Range("E40").Value = "  +0.69%  "
$ws.Range("D41").Value = "'0.05286"
$ws.Range("E41").Value = "  +0.15%  "
$ws.Range("D42").Value = "'6.723"
$ws.Range("D43").Value = "'116.81"
$ws.Range("E43").Value = "  +4.15%  "
$ws.Range("D44").Value = "'8.505"
$ws.Range("E44").Value = "  -0.32%  "
$ws.Range("D45").Value = "'0.1515"
$ws.Range("E45").Value = "  -0.21%  "
$ws.Range("D46").Value = "'0.4807"
$ws.Range("E46").Value = "  +0.35%  "
$ws.Range("D47").Value = "'10.46"
$ws.Range("E47").Value = "  -1.14%  "
$ws.Range("D48").Value = "'0.9996"
$ws.Range("E48").Value = "  -0.03%  "
$ws.Range("E49").Value = "  -1.20%  "
$ws.Range("D50").Value = "'66.71"
$ws.Range("E50").Value = "  -1.04%  "
$ws.Range("E51").Value = "  -0.90%  "
